# Add data for 2022-01-05
# - Rename sheet / update "through December 27" -> "through December 28"
# - Update December counts across years (2015-2021) for several neighborhoods
#   to reflect the newly added day of data (one extra carjacking date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet (tab) name and the header label in B1 (shared string) to
# reflect the new "through" date.
$ws.Name = "Through 2021-12-28"
$ws.Range("B1").Value = "December 2021 (through December 28)"

# Updated / newly populated monthly counts.
$ws.Range("AL3").Value  = 1    # Englewood            - December 2018
$ws.Range("BJ3").Value  = 3    # Englewood             - December 2016
$ws.Range("N4").Value   = 13   # North Lawndale        - December 2020
$ws.Range("AL7").Value  = 14   # Austin                - December 2018
$ws.Range("BV7").Value  = 6    # Austin                - December 2015
$ws.Range("BJ8").Value  = 5    # Chatham               - December 2016
$ws.Range("AX9").Value  = 3    # Grand Crossing        - December 2017
$ws.Range("N10").Value  = 3    # Douglas               - December 2020
$ws.Range("AX11").Value = 5    # Humboldt Park         - December 2017
$ws.Range("N17").Value  = 2    # United Center         - December 2020
$ws.Range("N21").Value  = 3    # Wicker Park           - December 2020
$ws.Range("AL23").Value = 4    # Little Village        - December 2018
$ws.Range("BJ24").Value = 6    # South Shore           - December 2016
$ws.Range("B26").Value  = 4    # Lake View             - December 2021 (through Dec 28)
$ws.Range("AL26").Value = 4    # Lake View             - December 2018
$ws.Range("Z29").Value  = 2    # Avalon Park           - December 2019
$ws.Range("B38").Value  = 3    # Auburn Gresham        - December 2021 (through Dec 28)
$ws.Range("N40").Value  = 5    # Calumet Heights       - December 2020
$ws.Range("B48").Value  = 1    # East Side             - December 2021 (through Dec 28)
$ws.Range("Z93").Value  = 1    # River North           - December 2019
$ws.Range("BV96").Value = 1    # South Deering         - December 2015
